$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Vorname" (first name) and "Nachname" (last name) columns were
# swapped: column B now holds what used to be in column C, and vice
# versa - including the header cells in row 1.
$lastRow = 181
for ($r = 1; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value2 = $cVal
    $ws.Cells.Item($r, 3).Value2 = $bVal
}

# Update the active selection to match the author's saved cursor position.
[void]$ws.Range("E4").Select()
